$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp header
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 05:22"

# Reorder: Mexico now ranks above Arabia Saudita and Filipinas (rows 37-39)
$ws.Range("A37").Value = "Mexico"
$ws.Range("B37").Value = 5399
$ws.Range("C37").Value = 385
$ws.Range("D37").Value = 2125
$ws.Range("E37").Value = 2868
$ws.Range("F37").Value = 207
$ws.Range("G37").Value = 74
$ws.Range("H37").Value = 406

$ws.Range("A38").Value = "Arabia Saudita"
$ws.Range("B38").Value = 5369
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 889
$ws.Range("E38").Value = 4407
$ws.Range("F38").Value = 59
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 73

$ws.Range("A39").Value = "Filipinas"
$ws.Range("B39").Value = 5223
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 295
$ws.Range("E39").Value = 4593
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 335

# Other straightforward numeric updates
$ws.Range("B17").Value = 25684
$ws.Range("C17").Value = 422
$ws.Range("E17").Value = 10106
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 1552

$ws.Range("B34").Value = 6416
$ws.Range("C34").Value = 16
$ws.Range("E34").Value = 2757

$ws.Range("B99").Value = 419
$ws.Range("C99").Value = 12
$ws.Range("E99").Value = 381
$ws.Range("G99").Value = 5
$ws.Range("H99").Value = 31

$ws.Range("B114").Value = 267
$ws.Range("C114").Value = 1
$ws.Range("E114").Value = 98

$ws.Range("B124").Value = 158
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 73
$ws.Range("E124").Value = 77
$ws.Range("F124").Value = 17
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 8

$ws.Range("D131").Value = 96
$ws.Range("E131").Value = 26

# Reorder: Zimbabue now ranks above Islas Virgenes de los Estados Unidos (rows 177-178)
$ws.Range("A177").Value = "Zimbabue"
$ws.Range("B177").Value = 18
$ws.Range("C177").Value = 1
$ws.Range("D177").Value = 1
$ws.Range("E177").Value = 14
$ws.Range("H177").Value = 3

$ws.Range("A178").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B178").Value = 17
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 17
$ws.Range("H178").Value = 0
